# Remove the "Validity Testing" comment left on "Properties" by
# Pontolillo, Gabriel J. (author-name scrubbing commit). Deleting the
# comment via the Comments collection removes the commentRangeStart/
# commentRangeEnd/commentReference markup from the body and drops the
# comment's own content from the comments part.
$d = $word.ActiveDocument

for ($i = $d.Comments.Count; $i -ge 1; $i--) {
    $d.Comments.Item($i).Delete()
}

# With the comment reference gone, "Properties" and ":" are two
# adjacent runs that used to be split only to host the comment anchor.
# Re-running Find/Replace across them collapses the text back into a
# single run, same as Word does when it coalesces runs around a
# deleted comment range.
$d.Content.Find.Execute("Properties:", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "Properties:", 2)
